# Regenerate catalog content only to point to new data products.
#
# The ShareURL column (N) previously held one ArcGIS dashboard link per
# county, each built from the same dashboard id with a "#geoName=..."
# fragment appended. The dashboards have been regenerated server-side, so
# each block of 11 rows (one per GeoName) now shares a single, plain
# dashboard URL (no fragment) - three new dashboard ids total, one per
# 11-row block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blockUrls = @(
    "https://www.arcgis.com/apps/dashboards/025b75f3490b4e79ae764e2c27c09a06",
    "https://www.arcgis.com/apps/dashboards/d1bb7ef5468f495788703352b1c5f896",
    "https://www.arcgis.com/apps/dashboards/2d3cc0a173d949f0a1a39146b37e1831"
)

$firstRow = 2
$lastRow = 34
$blockSize = 11

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $blockIndex = [Math]::Floor(($row - $firstRow) / $blockSize)
    $url = $blockUrls[$blockIndex]
    $ws.Cells.Item($row, 14).Value2 = $url
}
